$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook is a "handback status" report. A new handback run replaced
# two source-file GUIDs and regenerated their xliff hashes / timestamps:
#   74c554aa-24bd-440b-b828-92a8771f16fc  ->  1b32a180-1177-4f4a-b586-609f6897c251
#   ddd3338b-bca7-4a11-a407-98ddaabb3629  ->  ffffa128b606-038d-4129-905b-5be7893b88f8
# and the xliff content hash:
#   08dae925d280fb6cf46d30236e0fce1edfe32ac1 -> 3150de2df14ed3312291f43aba5fdad40b496d3e
# (the second row's handoff/handback xlf now coincides with the first row's,
# since its own old hash is gone from the regenerated report).
# ---------------------------------------------------------------------------

$guidNew1 = "1b32a180-1177-4f4a-b586-609f6897c251"
$guidNew2 = "ffffa128b606-038d-4129-905b-5be7893b88f8"

# ===== Sheet "Overview" =====
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$guidNew1.md"
$ws1.Range("B2").Value = "e2e\$guidNew1.md"
$ws1.Range("G2").Value = "2016-08-16 02:56:20"

$ws1.Range("A3").Value = "$guidNew2.md"
$ws1.Range("B3").Value = "e2e\$guidNew2.md"
$ws1.Range("G3").Value = "2016-08-16 02:56:20"

# Rebuild the hyperlinks so their visible "display" text reflects the new
# file names (the link targets themselves are untouched by this report run).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/660d17c1b3b91c685c7163541e028bfe61599898/e2e/74c554aa-24bd-440b-b828-92a8771f16fc.md", "", "", "e2e\$guidNew1.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/660d17c1b3b91c685c7163541e028bfe61599898/e2e/ddd3338b-bca7-4a11-a407-98ddaabb3629.md", "", "", "e2e\$guidNew2.md")

# ===== Sheet "zh-cn" =====
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$guidNew1.md"
$ws2.Range("G2").Value = "$guidNew1.3150de2df14ed3312291f43aba5fdad40b496d3e.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-16 02:56:15"
$ws2.Range("I2").Value = "$guidNew1.md"
$ws2.Range("J2").Value = "$guidNew1.3150de2df14ed3312291f43aba5fdad40b496d3e.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-16 02:56:32"

$ws2.Range("A3").Value = "$guidNew2.md"
$ws2.Range("G3").Value = "$guidNew1.3150de2df14ed3312291f43aba5fdad40b496d3e.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-16 02:56:15"
$ws2.Range("I3").Value = "$guidNew2.md"
$ws2.Range("J3").Value = "$guidNew1.3150de2df14ed3312291f43aba5fdad40b496d3e.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-16 02:56:32"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/660d17c1b3b91c685c7163541e028bfe61599898/e2e/74c554aa-24bd-440b-b828-92a8771f16fc.md", "", "", "$guidNew1.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0ae7d0431fe9fc0c4eed4222d6a468c53e98f5f5/e2e/74c554aa-24bd-440b-b828-92a8771f16fc.md", "", "", "$guidNew1.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/660d17c1b3b91c685c7163541e028bfe61599898/e2e/ddd3338b-bca7-4a11-a407-98ddaabb3629.md", "", "", "$guidNew2.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0ae7d0431fe9fc0c4eed4222d6a468c53e98f5f5/e2e/ddd3338b-bca7-4a11-a407-98ddaabb3629.md", "", "", "$guidNew2.md")

# ===== Sheet "de-de" =====
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$guidNew1.md"
$ws3.Range("G2").Value = "$guidNew1.3150de2df14ed3312291f43aba5fdad40b496d3e.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-16 02:56:20"
$ws3.Range("I2").Value = "$guidNew1.md"
$ws3.Range("J2").Value = "$guidNew1.3150de2df14ed3312291f43aba5fdad40b496d3e.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-16 02:56:39"

$ws3.Range("A3").Value = "$guidNew2.md"
$ws3.Range("G3").Value = "$guidNew1.3150de2df14ed3312291f43aba5fdad40b496d3e.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-16 02:56:20"
$ws3.Range("I3").Value = "$guidNew2.md"
$ws3.Range("J3").Value = "$guidNew1.3150de2df14ed3312291f43aba5fdad40b496d3e.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-16 02:56:39"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/660d17c1b3b91c685c7163541e028bfe61599898/e2e/74c554aa-24bd-440b-b828-92a8771f16fc.md", "", "", "$guidNew1.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/04f1608afa2a115b1d9216146943a4da3717abd9/e2e/74c554aa-24bd-440b-b828-92a8771f16fc.md", "", "", "$guidNew1.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/660d17c1b3b91c685c7163541e028bfe61599898/e2e/ddd3338b-bca7-4a11-a407-98ddaabb3629.md", "", "", "$guidNew2.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/04f1608afa2a115b1d9216146943a4da3717abd9/e2e/ddd3338b-bca7-4a11-a407-98ddaabb3629.md", "", "", "$guidNew2.md")
